# Update leve profit calculation figures across all Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 640
$ws.Range("I6").Value = 640
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1920
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1808
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
# Row 70
$ws.Range("H70").Value = 78122.08
$ws.Range("I70").Value = 334433.34
$ws.Range("J70").Value = 1228.7
$ws.Range("K70").Value = 1003300.02
$ws.Range("L70").Value = 3686.1
$ws.Range("M70").Value = -1003030.02
$ws.Range("N70").Value = -4226.1
# Row 73
$ws.Range("H73").Value = 78122.08
$ws.Range("I73").Value = 334433.34
$ws.Range("J73").Value = 1228.7
$ws.Range("K73").Value = 1003300.02
$ws.Range("L73").Value = 3686.1
$ws.Range("M73").Value = -1002364.02
$ws.Range("N73").Value = -5558.1
# Row 111
$ws.Range("H111").Value = 860
$ws.Range("I111").Value = 860
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2580
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 487

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 827.05
$ws.Range("I2").Value = 561.71875
$ws.Range("J2").Value = 1888.375
$ws.Range("K2").Value = 561.71875
$ws.Range("L2").Value = 1888.375
$ws.Range("M2").Value = -448.71875
$ws.Range("N2").Value = -2114.375
# Row 45
$ws.Range("H45").Value = 1934.6666
$ws.Range("I45").Value = 2062.4
$ws.Range("J45").Value = 1775
$ws.Range("K45").Value = 2062.4
$ws.Range("L45").Value = 1775
$ws.Range("M45").Value = -1685.4
$ws.Range("N45").Value = -2529
# Row 61
$ws.Range("H61").Value = 2287.75
$ws.Range("I61").Value = 1865.2222
$ws.Range("J61").Value = 3555.3333
$ws.Range("K61").Value = 1865.2222
$ws.Range("L61").Value = 3555.3333
$ws.Range("M61").Value = -1653.2222
$ws.Range("N61").Value = -3979.3333
# Row 62
$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
# Row 63
$ws.Range("H63").Value = 7750
$ws.Range("I63").Value = 5500
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 5500
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -4814
$ws.Range("N63").Value = -11372
# Row 65
$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240
# Row 66
$ws.Range("H66").Value = 7750
$ws.Range("I66").Value = 5500
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 27500
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -24068
$ws.Range("N66").Value = -56864
# Row 74
$ws.Range("H74").Value = 1417.2941
$ws.Range("I74").Value = 1238.0769
$ws.Range("J74").Value = 1999.75
$ws.Range("K74").Value = 1238.0769
$ws.Range("L74").Value = 1999.75
$ws.Range("M74").Value = -364.0769
$ws.Range("N74").Value = -3747.75
# Row 77
$ws.Range("H77").Value = 1417.2941
$ws.Range("I77").Value = 1238.0769
$ws.Range("J77").Value = 1999.75
$ws.Range("K77").Value = 6190.3845
$ws.Range("L77").Value = 9998.75
$ws.Range("M77").Value = -1822.3845
$ws.Range("N77").Value = -18734.75
# Row 97
$ws.Range("H97").Value = 621.7222
$ws.Range("I97").Value = 656.6667
$ws.Range("J97").Value = 551.8333
$ws.Range("K97").Value = 656.6667
$ws.Range("L97").Value = 551.8333
$ws.Range("M97").Value = -160.6667
$ws.Range("N97").Value = -1543.8333
# Row 116
$ws.Range("H116").Value = 827.05
$ws.Range("I116").Value = 561.71875
$ws.Range("J116").Value = 1888.375
$ws.Range("K116").Value = 561.71875
$ws.Range("L116").Value = 1888.375
$ws.Range("M116").Value = 1732.28125
$ws.Range("N116").Value = -6476.375
# Row 122
$ws.Range("H122").Value = 2193.6924
$ws.Range("I122").Value = 2105.7917
$ws.Range("J122").Value = 3248.5
$ws.Range("K122").Value = 6317.375100000001
$ws.Range("L122").Value = 9745.5
$ws.Range("M122").Value = -3867.375100000001
$ws.Range("N122").Value = -14645.5
# Row 132
$ws.Range("H132").Value = 871312.6
$ws.Range("I132").Value = 1334231.9
$ws.Range("J132").Value = 3339
$ws.Range("K132").Value = 4002695.7
$ws.Range("L132").Value = 10017
$ws.Range("M132").Value = -4000165.7
$ws.Range("N132").Value = -15077
# Row 136
$ws.Range("H136").Value = 2287.75
$ws.Range("I136").Value = 1865.2222
$ws.Range("J136").Value = 3555.3333
$ws.Range("K136").Value = 5595.6666
$ws.Range("L136").Value = 10665.9999
$ws.Range("M136").Value = -3045.6666
$ws.Range("N136").Value = -15765.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 827.05
$ws.Range("I3").Value = 561.71875
$ws.Range("J3").Value = 1888.375
$ws.Range("K3").Value = 561.71875
$ws.Range("L3").Value = 1888.375
$ws.Range("M3").Value = -447.71875
$ws.Range("N3").Value = -2116.375
# Row 26
$ws.Range("H26").Value = 500000
$ws.Range("I26").Value = 500000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 500000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -499708
# Row 96
$ws.Range("H96").Value = 335000
$ws.Range("I96").Value = 335000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 335000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -332254
# Row 134
$ws.Range("H134").Value = 447288.4
$ws.Range("I134").Value = 743094.9399999999
$ws.Range("J134").Value = 3578.5557
$ws.Range("K134").Value = 2229284.82
$ws.Range("L134").Value = 10735.6671
$ws.Range("M134").Value = -2226749.82
$ws.Range("N134").Value = -15805.6671

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2946.4355
$ws.Range("I31").Value = 2267.4412
$ws.Range("J31").Value = 3770.9285
$ws.Range("K31").Value = 2267.4412
$ws.Range("L31").Value = 3770.9285
$ws.Range("M31").Value = -1972.4412
$ws.Range("N31").Value = -4360.9285
# Row 34
$ws.Range("H34").Value = 2946.4355
$ws.Range("I34").Value = 2267.4412
$ws.Range("J34").Value = 3770.9285
$ws.Range("K34").Value = 2267.4412
$ws.Range("L34").Value = 3770.9285
$ws.Range("M34").Value = -2065.4412
$ws.Range("N34").Value = -4174.9285
# Row 70
$ws.Range("H70").Value = 11960
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 11960
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 11960
$ws.Range("N70").Value = -12590
# Row 73
$ws.Range("H73").Value = 11960
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 11960
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 11960
$ws.Range("N73").Value = -14144

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1940.4584
$ws.Range("I5").Value = 2847.5833
$ws.Range("J5").Value = 1033.3334
$ws.Range("K5").Value = 8542.749899999999
$ws.Range("L5").Value = 3100.0002
$ws.Range("M5").Value = -8430.749899999999
$ws.Range("N5").Value = -3324.0002
# Row 7
$ws.Range("H7").Value = 105.85714
$ws.Range("I7").Value = 106.833336
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 320.500008
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -208.500008
$ws.Range("N7").Value = -524
# Row 92
$ws.Range("H92").Value = 1000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 3000
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5496
# Row 132
$ws.Range("H132").Value = 1808.303
$ws.Range("I132").Value = 980.6667
$ws.Range("J132").Value = 2281.238
$ws.Range("K132").Value = 8826.0003
$ws.Range("L132").Value = 20531.142
$ws.Range("M132").Value = -6296.0003
$ws.Range("N132").Value = -25591.142
# Row 133
$ws.Range("H133").Value = 3555.2
$ws.Range("I133").Value = 1932
$ws.Range("J133").Value = 5990
$ws.Range("K133").Value = 5796
$ws.Range("L133").Value = 17970
$ws.Range("M133").Value = -736
$ws.Range("N133").Value = -28090
# Row 134
$ws.Range("H134").Value = 3873.647
$ws.Range("I134").Value = 1824.1364
$ws.Range("J134").Value = 7631.0835
$ws.Range("K134").Value = 5472.4092
$ws.Range("L134").Value = 22893.2505
$ws.Range("M134").Value = -402.4092000000001
$ws.Range("N134").Value = -33033.25049999999
# Row 135
$ws.Range("H135").Value = 1940.4584
$ws.Range("I135").Value = 2847.5833
$ws.Range("J135").Value = 1033.3334
$ws.Range("K135").Value = 25628.2497
$ws.Range("L135").Value = 9300.000599999999
$ws.Range("M135").Value = -23093.2497
$ws.Range("N135").Value = -14370.0006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 109
$ws.Range("H109").Value = 22141.428
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 22141.428
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 22141.428
$ws.Range("N109").Value = -24221.428
# Row 122
$ws.Range("H122").Value = 3487.2258
$ws.Range("I122").Value = 3635.1538
$ws.Range("J122").Value = 2718
$ws.Range("K122").Value = 10905.4614
$ws.Range("L122").Value = 8154
$ws.Range("M122").Value = -8455.4614
$ws.Range("N122").Value = -13054
# Row 123
$ws.Range("H123").Value = 20377.46
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20377.46
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20377.46
$ws.Range("N123").Value = -25277.46
# Row 132
$ws.Range("H132").Value = 2208.558
$ws.Range("I132").Value = 1975.3572
$ws.Range("J132").Value = 2643.8667
$ws.Range("K132").Value = 5926.071599999999
$ws.Range("L132").Value = 7931.6001
$ws.Range("M132").Value = -3396.071599999999
$ws.Range("N132").Value = -12991.6001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 92
$ws.Range("H92").Value = 32000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 32000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 32000
$ws.Range("N92").Value = -36992
# Row 132
$ws.Range("H132").Value = 4835.927
$ws.Range("I132").Value = 5343.926
$ws.Range("J132").Value = 3856.2144
$ws.Range("K132").Value = 16031.778
$ws.Range("L132").Value = 11568.6432
$ws.Range("M132").Value = -13501.778
$ws.Range("N132").Value = -16628.6432
# Row 136
$ws.Range("H136").Value = 3194.0566
$ws.Range("I136").Value = 3239.0625
$ws.Range("J136").Value = 2762
$ws.Range("K136").Value = 9717.1875
$ws.Range("L136").Value = 8286
$ws.Range("M136").Value = -7167.1875
$ws.Range("N136").Value = -13386
# Row 138
$ws.Range("H138").Value = 87770
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 87770
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 87770
$ws.Range("N138").Value = -98050

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 3500
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -2888
$ws.Range("N2").Value = -4224
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 127
$ws.Range("H127").Value = 17770
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 17770
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 17770
$ws.Range("N127").Value = -27690
# Row 132
$ws.Range("H132").Value = 2323.1667
$ws.Range("I132").Value = 1761.4445
$ws.Range("J132").Value = 3165.75
$ws.Range("K132").Value = 5284.333500000001
$ws.Range("L132").Value = 9497.25
$ws.Range("M132").Value = -2754.333500000001
$ws.Range("N132").Value = -14557.25

